$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the account-summary figures ---------------------------------
# "Valor Mora" (overdue amount) grew because a new overdue period was added
$ws.Range("E11").Value = 140480
# "Cant. Periodos" (number of overdue periods) went from 3 to 4
$ws.Range("F13").Value = 4

# --- Make room for the new period row ------------------------------------
# Rows 19-22 are already blank, so inserting a row right above the
# signature block (row 23) shifts it down to rows 24-25 without disturbing
# the existing formatting of the data rows (16-18).
$ws.Rows("23").Insert()

# --- Add the new period (part 1 of the new account statement) -----------
# Row 19 becomes the new "last" data row: clone row 18's content/format.
$ws.Range("B18:J18").Copy($ws.Range("B19:J19"))

# Row 18 is no longer the last row, so it takes on the regular row style
# (same formatting as rows 16 and 17).
$ws.Range("B16:J16").Copy($ws.Range("B18:J18"))

# --- Resequence the period labels: 2507,2506,2505 -> 2505,2506,2507,2508
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"
